$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force text format on the Price/Volume columns so that
# numeric-looking strings (e.g. "1.001") are stored as text, matching
# the original inline-string cell type, instead of being parsed as numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '29.224.01'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '1.842.22'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '242.60'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').Value = '0.6636'
$ws.Range('E6').Value = '  -1.10%  '
$ws.Range('D8').Value = '44.80'
$ws.Range('D9').Value = '0.07452'
$ws.Range('E9').Value = '  +0.50%  '
$ws.Range('D11').Value = '23.40'
$ws.Range('E11').Value = '  +1.91%  '
$ws.Range('D12').Value = '0.07773'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('D13').Value = '1.856.23'
$ws.Range('E13').Value = '  +2.70%  '
$ws.Range('D14').Value = '5.025'
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').Value = '0.6729'
$ws.Range('E15').Value = '  -0.99%  '
$ws.Range('D16').Value = '83.54'
$ws.Range('E16').Value = '  -3.39%  '
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('D18').Value = '0.000008619'
$ws.Range('E18').Value = '  +4.22%  '
$ws.Range('D19').Value = '29.244.60'
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').Value = '2.112.31'
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('D21').Value = '227.45'
$ws.Range('E21').Value = '  -0.88%  '
$ws.Range('D22').Value = '12.56'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').Value = '7.194'
$ws.Range('E24').Value = '  -1.40%  '
$ws.Range('D25').Value = '1.001'
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').Value = '158.94'
$ws.Range('E26').Value = '  -0.86%  '
$ws.Range('D27').Value = '0.1408'
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('D28').Value = '8.639'
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('D31').Value = '4.137'
$ws.Range('E31').Value = '  -1.81%  '
$ws.Range('D32').Value = '4.061'
$ws.Range('E32').Value = '  -0.73%  '
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('D34').Value = '0.05337'
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('D35').Value = '1.877'
$ws.Range('E35').Value = '  +0.19%  '
$ws.Range('D36').Value = '0.7470'
$ws.Range('E36').Value = '  -0.76%  '
$ws.Range('D37').Value = '1.160'
$ws.Range('E37').Value = '  +1.80%  '
$ws.Range('D38').Value = '2.653'
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('D39').Value = '1.322.74'
$ws.Range('E39').Value = '  -0.75%  '
$ws.Range('D40').Value = '0.01803'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('D42').Value = '6.401'
$ws.Range('E42').Value = '  +7.04%  '
$ws.Range('D43').Value = '0.9075'
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('D44').Value = '0.9999'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').Value = '103.19'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '2.006.13'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').Value = '65.69'
$ws.Range('E47').Value = '  +2.72%  '
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('D50').Value = '0.07671'
$ws.Range('E50').Value = '  -8.21%  '
$ws.Range('D51').Value = '1.755'
$ws.Range('E51').Value = '  -0.55%  '

# Restore the default (Normal) style so the cells keep their original,
# unstyled appearance once the text values have been committed.
$dataRange.Style = "Normal"
